$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header value tweaks
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 ("CON") - remove the now-unused Lichtwark values, keep/update C2
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = 26.188281108599256
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()

# Row 3 ("STR") - updated passive values
$ws.Range("B3").Value = 23.27312537304482
$ws.Range("C3").Value = 34.209238566352383
$ws.Range("D3").Value = 35.843641919456275
$ws.Range("E3").Value = 19.785949458981804

# Update selection to reflect the newly relevant range
$ws.Range("B1:E3").Select()
